# Updated cryptos list values (price / 1h volume change) per scraped diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> @{ D = newPriceText (optional); E = newVolumeText }
$updates = @{
    2  = @{ D = "69.041.43";  E = "  +2.02%  " }
    3  = @{ D = "3.932.20";   E = "  +1.11%  " }
    4  = @{           E = "  +0.23%  " }
    5  = @{ D = "484.72";     E = "  +3.90%  " }
    6  = @{ D = "147.56";     E = "  -0.75%  " }
    7  = @{           E = "  -2.17%  " }
    8  = @{           E = "  +0.12%  " }
    9  = @{ D = "0.724";      E = "  -3.20%  " }
    10 = @{ D = "0.172";      E = "  +9.54%  " }
    11 = @{           E = "  +13.07%  " }
    12 = @{ D = "42.64";      E = "  -2.71%  " }
    13 = @{           E = "  +1.02%  " }
    14 = @{ D = "4.568.06";   E = "  +1.11%  " }
    15 = @{ D = "3.955.68";   E = "  +2.48%  " }
    16 = @{ D = "14.55";      E = "  -1.70%  " }
    17 = @{           E = "  -0.30%  " }
    18 = @{ D = "19.70";      E = "  -1.92%  " }
    19 = @{           E = "  -3.09%  " }
    20 = @{ D = "69.166.13";  E = "  +2.00%  " }
    21 = @{ D = "435.57";     E = "  +1.00%  " }
    22 = @{           E = "  -1.79%  " }
    23 = @{           E = "  +0.62%  " }
    24 = @{ D = "87.44";      E = "  -1.28%  " }
    25 = @{ D = "11.61";      E = "  +14.55%  " }
    26 = @{           E = "  -0.60%  " }
    27 = @{ D = "10.59";      E = "  +2.67%  " }
    28 = @{ D = "38.16";      E = "  +0.64%  " }
    29 = @{ D = "5.88";       E = "  +6.88%  " }
    30 = @{ D = "714.45";     E = "  -2.49%  " }
    31 = @{ D = "13.23";      E = "  -4.11%  " }
    32 = @{           E = "  -4.95%  " }
    33 = @{           E = "  +2.60%  " }
    34 = @{ D = "0.0₃0905";   E = "  +32.54%  " }
    35 = @{ D = "41.18";      E = "  -4.42%  " }
    36 = @{ D = "58.62";      E = "  +1.23%  " }
    37 = @{ D = "0.151";      E = "  -6.63%  " }
    38 = @{           E = "  +1.10%  " }
    39 = @{           E = "  -0.18%  " }
    40 = @{ D = "0.0471" }
    41 = @{ D = "2.77";       E = "  +7.35%  " }
    42 = @{ D = "2.99";       E = "  +7.19%  " }
    43 = @{           E = "  +1.64%  " }
    44 = @{ D = "0.338";      E = "  -2.13%  " }
    45 = @{           E = "  -1.32%  " }
    46 = @{           E = "  +0.18%  " }
    47 = @{ D = "3.41";       E = "  -0.91%  " }
    48 = @{ D = "2.15";       E = "  +0.40%  " }
    49 = @{ D = "147.88";     E = "  +2.36%  " }
    50 = @{ D = "3.15";       E = "  -2.88%  " }
    51 = @{           E = "  -2.37%  " }
}

foreach ($row in $updates.Keys) {
    $cellUpdates = $updates[$row]
    if ($cellUpdates.ContainsKey("D")) {
        $ws.Range("D$row").Value = $cellUpdates["D"]
    }
    if ($cellUpdates.ContainsKey("E")) {
        $ws.Range("E$row").Value = $cellUpdates["E"]
    }
}
